$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2450980392156863
$ws.Range("C2").Value = 0.4248366013071895
$ws.Range("J2").Value = 0.03594771241830065
$ws.Range("P2").Value = 0.1699346405228758
$ws.Range("S2").Value = 0.1241830065359477
$ws.Range("B3").Value = 0.0218978102189781
$ws.Range("C3").Value = 0.0291970802919708
$ws.Range("J3").Value = 0.05109489051094891
$ws.Range("P3").Value = 0.6788321167883211
$ws.Range("S3").Value = 0.218978102189781
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.6041666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.07112970711297072
$ws.Range("D6").Value = 0.008368200836820083
$ws.Range("F6").Value = 0.08786610878661087
$ws.Range("J6").Value = 0.3096234309623431
$ws.Range("O6").Value = 0.01673640167364017
$ws.Range("Q6").Value = 0.1464435146443515
$ws.Range("R6").Value = 0.05857740585774059
$ws.Range("S6").Value = 0.301255230125523
$ws.Range("B7").Value = 0.1269035532994924
$ws.Range("D7").Value = 0.005076142131979695
$ws.Range("F7").Value = 0.08121827411167512
$ws.Range("J7").Value = 0.1116751269035533
$ws.Range("O7").Value = 0.005076142131979695
$ws.Range("Q7").Value = 0.1319796954314721
$ws.Range("R7").Value = 0.08629441624365482
$ws.Range("S7").Value = 0.4517766497461929
$ws.Range("B8").Value = 0.08108108108108109
$ws.Range("D8").Value = 0.02252252252252252
$ws.Range("E8").Value = 0.002252252252252252
$ws.Range("F8").Value = 0.06981981981981981
$ws.Range("J8").Value = 0.1509009009009009
$ws.Range("O8").Value = 0.02477477477477477
$ws.Range("Q8").Value = 0.1756756756756757
$ws.Range("R8").Value = 0.08333333333333333
$ws.Range("S8").Value = 0.3896396396396397
$ws.Range("B9").Value = 0.07389162561576355
$ws.Range("D9").Value = 0.03448275862068965
$ws.Range("F9").Value = 0.04926108374384237
$ws.Range("J9").Value = 0.09852216748768473
$ws.Range("O9").Value = 0.03448275862068965
$ws.Range("Q9").Value = 0.1822660098522167
$ws.Range("R9").Value = 0.1133004926108374
$ws.Range("S9").Value = 0.4137931034482759
$ws.Range("B10").Value = 0.1038251366120219
$ws.Range("D10").Value = 0.02029664324746292
$ws.Range("E10").Value = 0.00156128024980484
$ws.Range("F10").Value = 0.06713505074160812
$ws.Range("J10").Value = 0.1366120218579235
$ws.Range("O10").Value = 0.01717408274785324
$ws.Range("Q10").Value = 0.1912568306010929
$ws.Range("R10").Value = 0.09289617486338798
$ws.Range("S10").Value = 0.3692427790788447
$ws.Range("G11").Value = 0.1335616438356164
$ws.Range("J11").Value = 0.08904109589041095
$ws.Range("K11").Value = 0.1883561643835616
$ws.Range("L11").Value = 0.5684931506849316
$ws.Range("S11").Value = 0.02054794520547945
$ws.Range("G12").Value = 0.7283236994219653
$ws.Range("J12").Value = 0.2023121387283237
$ws.Range("K12").Value = 0.0115606936416185
$ws.Range("L12").Value = 0.03468208092485549
$ws.Range("S12").Value = 0.02312138728323699
$ws.Range("G13").Value = 0.76
$ws.Range("J13").Value = 0.08
$ws.Range("S13").Value = 0.16
$ws.Range("F15").Value = 0.02586206896551724
$ws.Range("H15").Value = 0.1551724137931035
$ws.Range("I15").Value = 0.09482758620689655
$ws.Range("J15").Value = 0.3103448275862069
$ws.Range("K15").Value = 0.06896551724137931
$ws.Range("M15").Value = 0.01293103448275862
$ws.Range("O15").Value = 0.08189655172413793
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.01764705882352941
$ws.Range("H16").Value = 0.1
$ws.Range("I16").Value = 0.05294117647058823
$ws.Range("J16").Value = 0.488235294117647
$ws.Range("K16").Value = 0.05294117647058823
$ws.Range("M16").Value = 0.04705882352941176
$ws.Range("O16").Value = 0.03529411764705882
$ws.Range("S16").Value = 0.2058823529411765
$ws.Range("F17").Value = 0.02107728337236534
$ws.Range("H17").Value = 0.1850117096018735
$ws.Range("I17").Value = 0.1124121779859485
$ws.Range("J17").Value = 0.3864168618266979
$ws.Range("K17").Value = 0.1053864168618267
$ws.Range("M17").Value = 0.01873536299765808
$ws.Range("O17").Value = 0.06088992974238876
$ws.Range("S17").Value = 0.1100702576112412
$ws.Range("F18").Value = 0.02392344497607655
$ws.Range("H18").Value = 0.1818181818181818
$ws.Range("I18").Value = 0.1052631578947368
$ws.Range("J18").Value = 0.4019138755980861
$ws.Range("K18").Value = 0.1100478468899522
$ws.Range("M18").Value = 0.02392344497607655
$ws.Range("O18").Value = 0.03827751196172249
$ws.Range("S18").Value = 0.1148325358851675
$ws.Range("F19").Value = 0.02380952380952381
$ws.Range("H19").Value = 0.2073732718894009
$ws.Range("I19").Value = 0.07910906298003072
$ws.Range("J19").Value = 0.347926267281106
$ws.Range("K19").Value = 0.1098310291858679
$ws.Range("M19").Value = 0.02150537634408602
$ws.Range("O19").Value = 0.07450076804915515
$ws.Range("S19").Value = 0.1359447004608295
